# Update Sheets per scheduled runner refresh of market-board data
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 45342.668
$ws.Range("I62").Value = 74119.21000000001
$ws.Range("J62").Value = 5055.5
$ws.Range("K62").Value = 74119.21000000001
$ws.Range("L62").Value = 5055.5
$ws.Range("M62").Value = -73495.21000000001
$ws.Range("N62").Value = -6303.5
$ws.Range("H65").Value = 45342.668
$ws.Range("I65").Value = 74119.21000000001
$ws.Range("J65").Value = 5055.5
$ws.Range("K65").Value = 370596.05
$ws.Range("L65").Value = 25277.5
$ws.Range("M65").Value = -367476.05
$ws.Range("N65").Value = -31517.5
$ws.Range("H70").Value = 1814.7273
$ws.Range("I70").Value = 1260.3334
$ws.Range("J70").Value = 2480
$ws.Range("K70").Value = 3781.0002
$ws.Range("L70").Value = 7440
$ws.Range("M70").Value = -3511.0002
$ws.Range("N70").Value = -7980
$ws.Range("H73").Value = 1814.7273
$ws.Range("I73").Value = 1260.3334
$ws.Range("J73").Value = 2480
$ws.Range("K73").Value = 3781.0002
$ws.Range("L73").Value = 7440
$ws.Range("M73").Value = -2845.0002
$ws.Range("N73").Value = -9312
$ws.Range("H74").Value = 4330177
$ws.Range("I74").Value = 5194812.5
$ws.Range("K74").Value = 5194812.5
$ws.Range("M74").Value = -5193876.5
$ws.Range("H77").Value = 4330177
$ws.Range("I77").Value = 5194812.5
$ws.Range("K77").Value = 25974062.5
$ws.Range("M77").Value = -25969382.5
$ws.Range("H82").Value = 40004116
$ws.Range("I82").Value = 2847.3333
$ws.Range("J82").Value = 100006024
$ws.Range("K82").Value = 8541.999899999999
$ws.Range("L82").Value = 300018072
$ws.Range("M82").Value = -8135.999899999999
$ws.Range("N82").Value = -300018884
$ws.Range("H85").Value = 40004116
$ws.Range("I85").Value = 2847.3333
$ws.Range("J85").Value = 100006024
$ws.Range("K85").Value = 8541.999899999999
$ws.Range("L85").Value = 300018072
$ws.Range("M85").Value = -7137.999899999999
$ws.Range("N85").Value = -300020880
$ws.Range("H92").Value = 833
$ws.Range("I92").Value = 768.8
$ws.Range("K92").Value = 768.8
$ws.Range("M92").Value = 479.2
$ws.Range("H107").Value = 909215.75
$ws.Range("I107").Value = 1000107.3
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 1000107.3
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = -998187.3
$ws.Range("N107").Value = -4140
$ws.Range("H116").Value = 3820.8462
$ws.Range("I116").Value = 4345.625
$ws.Range("K116").Value = 4345.625
$ws.Range("M116").Value = -903.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1368.4
$ws.Range("I122").Value = 1170.5714
$ws.Range("J122").Value = 1830
$ws.Range("K122").Value = 3511.7142
$ws.Range("L122").Value = 5490
$ws.Range("M122").Value = -1061.7142
$ws.Range("N122").Value = -10390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2936.3635
$ws.Range("I86").Value = 2512.5
$ws.Range("J86").Value = 4066.6667
$ws.Range("K86").Value = 2512.5
$ws.Range("L86").Value = 4066.6667
$ws.Range("M86").Value = -1389.5
$ws.Range("N86").Value = -6312.6667
$ws.Range("H89").Value = 2936.3635
$ws.Range("I89").Value = 2512.5
$ws.Range("J89").Value = 4066.6667
$ws.Range("K89").Value = 12562.5
$ws.Range("L89").Value = 20333.3335
$ws.Range("M89").Value = -6946.5
$ws.Range("N89").Value = -31565.3335
$ws.Range("H105").Value = 5201.15
$ws.Range("I105").Value = 4975.8125
$ws.Range("J105").Value = 6102.5
$ws.Range("K105").Value = 4975.8125
$ws.Range("L105").Value = 6102.5
$ws.Range("M105").Value = -3228.8125
$ws.Range("N105").Value = -9596.5
$ws.Range("H107").Value = 9602.3125
$ws.Range("I107").Value = 1389.1333
$ws.Range("K107").Value = 1389.1333
$ws.Range("M107").Value = 530.8667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1742.5834
$ws.Range("I16").Value = 1741.1
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 1741.1
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -1454.1
$ws.Range("N16").Value = -2324
$ws.Range("H113").Value = 1742.5834
$ws.Range("I113").Value = 1741.1
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1741.1
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 428.9000000000001
$ws.Range("N113").Value = -6090
$ws.Range("H134").Value = 2869.4348
$ws.Range("I134").Value = 1905.1052
$ws.Range("J134").Value = 7450
$ws.Range("K134").Value = 5715.3156
$ws.Range("L134").Value = 22350
$ws.Range("M134").Value = -3180.3156
$ws.Range("N134").Value = -27420
$ws.Range("H140").Value = 52137.777
$ws.Range("J140").Value = 52137.777
$ws.Range("L140").Value = 52137.777
$ws.Range("N140").Value = -62497.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 6561.375
$ws.Range("I80").Value = 3489
$ws.Range("J80").Value = 7000.2856
$ws.Range("K80").Value = 10467
$ws.Range("L80").Value = 21000.8568
$ws.Range("M80").Value = -9531
$ws.Range("N80").Value = -22872.8568
$ws.Range("H83").Value = 6561.375
$ws.Range("I83").Value = 3489
$ws.Range("J83").Value = 7000.2856
$ws.Range("K83").Value = 31401
$ws.Range("L83").Value = 63002.5704
$ws.Range("M83").Value = -26721
$ws.Range("N83").Value = -72362.5704
$ws.Range("H117").Value = 2349.0667
$ws.Range("I117").Value = 1252.6666
$ws.Range("J117").Value = 2623.1667
$ws.Range("K117").Value = 3757.9998
$ws.Range("L117").Value = 7869.500100000001
$ws.Range("M117").Value = -315.9998000000001
$ws.Range("N117").Value = -14753.5001
$ws.Range("H129").Value = 1955.5358
$ws.Range("I129").Value = 963.625
$ws.Range("J129").Value = 2352.3
$ws.Range("K129").Value = 2890.875
$ws.Range("L129").Value = 7056.900000000001
$ws.Range("M129").Value = 2109.125
$ws.Range("N129").Value = -17056.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4343.3335
$ws.Range("I70").Value = 4100
$ws.Range("J70").Value = 7750
$ws.Range("K70").Value = 4100
$ws.Range("L70").Value = 7750
$ws.Range("M70").Value = -3830
$ws.Range("N70").Value = -8290
$ws.Range("H73").Value = 4343.3335
$ws.Range("I73").Value = 4100
$ws.Range("J73").Value = 7750
$ws.Range("K73").Value = 4100
$ws.Range("L73").Value = 7750
$ws.Range("M73").Value = -3164
$ws.Range("N73").Value = -9622
$ws.Range("H116").Value = 26000
$ws.Range("J116").Value = 26000
$ws.Range("L116").Value = 26000
$ws.Range("N116").Value = -35178
$ws.Range("H122").Value = 1317769.1
$ws.Range("I122").Value = 2633258.8
$ws.Range("J122").Value = 2279.4
$ws.Range("K122").Value = 7899776.399999999
$ws.Range("L122").Value = 6838.200000000001
$ws.Range("M122").Value = -7897326.399999999
$ws.Range("N122").Value = -11738.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13500
$ws.Range("J41").Value = 13500
$ws.Range("L41").Value = 13500
$ws.Range("N41").Value = -14280
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H132").Value = 1195.3721
$ws.Range("I132").Value = 919.9429
$ws.Range("J132").Value = 2759.8287
$ws.Range("K132").Value = 2759.8287
$ws.Range("M132").Value = -229.8287
